$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DUA")

# Insert a new header row at the top of the DUA sheet, shifting the
# existing content (and the merged B-column dropdown range) down by one.
$ws.Rows.Item(1).Insert()

# New header row: "Field" / "Content" labels, bold red text (style used
# elsewhere in the workbook for section headers).
$ws.Range("A1").Value = "Field"
$ws.Range("B1").Value = "Content"
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").Font.Color = 255

# The row-insert shifts the merged range / dropdown cell down automatically,
# but the validation list formula needs to be repointed at the now-shifted
# D2:D6 source range.
[void]$ws.Range("B2").Validation.Modify(3, 1, 1, "=`$D`$2:`$D`$6")

# Match the selection left on the sheet and make DUA the active tab.
[void]$ws.Range("B3:B20").Select()
[void]$ws.Activate()
